$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 2054.1667
$ws_ALC.Range("J17").Value = 2025
$ws_ALC.Range("L17").Value = 6075
$ws_ALC.Range("N17").Value = -6411

# ALC row 29
$ws_ALC.Range("H29").Value = 1525
$ws_ALC.Range("I29").Value = 100
$ws_ALC.Range("K29").Value = 300
$ws_ALC.Range("M29").Value = -19

# ALC row 38
$ws_ALC.Range("H38").Value = 909.1111
$ws_ALC.Range("I38").Value = 36.4
$ws_ALC.Range("J38").Value = 2000
$ws_ALC.Range("K38").Value = 109.2
$ws_ALC.Range("L38").Value = 6000
$ws_ALC.Range("M38").Value = 262.8
$ws_ALC.Range("N38").Value = -6744

# ALC row 39
$ws_ALC.Range("H39").Value = 2037.6
$ws_ALC.Range("I39").Value = 47
$ws_ALC.Range("J39").Value = 10000
$ws_ALC.Range("K39").Value = 141
$ws_ALC.Range("L39").Value = 30000
$ws_ALC.Range("M39").Value = 155
$ws_ALC.Range("N39").Value = -30592

# ALC row 43
$ws_ALC.Range("I43").Value = 1994.5
$ws_ALC.Range("J43").Value = 1999.5
$ws_ALC.Range("K43").Value = 1994.5
$ws_ALC.Range("L43").Value = 1999.5
$ws_ALC.Range("M43").Value = -1925.5
$ws_ALC.Range("N43").Value = -2137.5

# ALC row 45
$ws_ALC.Range("H45").Value = 0
$ws_ALC.Range("J45").Value = 0
$ws_ALC.Range("L45").Value = 0
$ws_ALC.Range("N45").ClearContents()

# ALC row 86
$ws_ALC.Range("H86").Value = 3008.5
$ws_ALC.Range("I86").Value = 1701.5
$ws_ALC.Range("K86").Value = 1701.5
$ws_ALC.Range("M86").Value = -578.5

# ALC row 89
$ws_ALC.Range("H89").Value = 3008.5
$ws_ALC.Range("I89").Value = 1701.5
$ws_ALC.Range("K89").Value = 8507.5
$ws_ALC.Range("M89").Value = -2891.5

# ALC row 100
$ws_ALC.Range("H100").Value = 1000
$ws_ALC.Range("I100").Value = 1000
$ws_ALC.Range("K100").Value = 1000
$ws_ALC.Range("M100").Value = -459

# ALC row 137
$ws_ALC.Range("H137").Value = 1732.7693
$ws_ALC.Range("J137").Value = 1816.1666
$ws_ALC.Range("L137").Value = 5448.4998
$ws_ALC.Range("N137").Value = -10548.4998

# ARM row 61
$ws_ARM.Range("H61").Value = 4437.3335
$ws_ARM.Range("I61").Value = 3926.2
$ws_ARM.Range("J61").Value = 6993
$ws_ARM.Range("K61").Value = 3926.2
$ws_ARM.Range("L61").Value = 6993
$ws_ARM.Range("M61").Value = -3714.2
$ws_ARM.Range("N61").Value = -7417

# ARM row 74
$ws_ARM.Range("H74").Value = 6000
$ws_ARM.Range("I74").Value = 6000
$ws_ARM.Range("K74").Value = 6000
$ws_ARM.Range("M74").Value = -5126

# ARM row 77
$ws_ARM.Range("H77").Value = 6000
$ws_ARM.Range("I77").Value = 6000
$ws_ARM.Range("K77").Value = 30000
$ws_ARM.Range("M77").Value = -25632

# ARM row 136
$ws_ARM.Range("H136").Value = 4437.3335
$ws_ARM.Range("I136").Value = 3926.2
$ws_ARM.Range("J136").Value = 6993
$ws_ARM.Range("K136").Value = 11778.6
$ws_ARM.Range("L136").Value = 20979
$ws_ARM.Range("M136").Value = -9228.599999999999
$ws_ARM.Range("N136").Value = -26079

# BSM row 12
$ws_BSM.Range("H12").Value = 3590.1428
$ws_BSM.Range("I12").Value = 3336
$ws_BSM.Range("J12").Value = 3780.75
$ws_BSM.Range("K12").Value = 3336
$ws_BSM.Range("L12").Value = 3780.75
$ws_BSM.Range("M12").Value = -3168
$ws_BSM.Range("N12").Value = -4116.75

# BSM row 134
$ws_BSM.Range("H134").Value = 12332
$ws_BSM.Range("J134").Value = 11798.4
$ws_BSM.Range("L134").Value = 35395.2
$ws_BSM.Range("N134").Value = -40465.2

# CRP row 31
$ws_CRP.Range("H31").Value = 6971.6665
$ws_CRP.Range("I31").Value = 9056.666999999999
$ws_CRP.Range("J31").Value = 6276.6665
$ws_CRP.Range("K31").Value = 9056.666999999999
$ws_CRP.Range("L31").Value = 6276.6665
$ws_CRP.Range("M31").Value = -8761.666999999999
$ws_CRP.Range("N31").Value = -6866.6665

# CRP row 34
$ws_CRP.Range("H34").Value = 6971.6665
$ws_CRP.Range("I34").Value = 9056.666999999999
$ws_CRP.Range("J34").Value = 6276.6665
$ws_CRP.Range("K34").Value = 9056.666999999999
$ws_CRP.Range("L34").Value = 6276.6665
$ws_CRP.Range("M34").Value = -8854.666999999999
$ws_CRP.Range("N34").Value = -6680.6665

# CRP row 58
$ws_CRP.Range("H58").Value = 0
$ws_CRP.Range("I58").Value = 0
$ws_CRP.Range("K58").Value = 0
$ws_CRP.Range("M58").ClearContents()

# CRP row 132
$ws_CRP.Range("H132").Value = 0
$ws_CRP.Range("I132").Value = 0
$ws_CRP.Range("K132").Value = 0
$ws_CRP.Range("M132").ClearContents()

# CRP row 134
$ws_CRP.Range("H134").Value = 3969.4
$ws_CRP.Range("J134").Value = 4724.25
$ws_CRP.Range("L134").Value = 14172.75
$ws_CRP.Range("N134").Value = -19242.75

# CRP row 136
$ws_CRP.Range("H136").Value = 0
$ws_CRP.Range("I136").Value = 0
$ws_CRP.Range("K136").Value = 0
$ws_CRP.Range("M136").ClearContents()

# CUL row 122
$ws_CUL.Range("H122").Value = 2151
$ws_CUL.Range("I122").Value = 1438.75
$ws_CUL.Range("K122").Value = 12948.75
$ws_CUL.Range("M122").Value = -10498.75

# LTW row 40
$ws_LTW.Range("H40").Value = 5263.4375
$ws_LTW.Range("I40").Value = 4862.6924
$ws_LTW.Range("K40").Value = 4862.6924
$ws_LTW.Range("M40").Value = -4726.6924

# LTW row 68
$ws_LTW.Range("H68").Value = 2257.8572
$ws_LTW.Range("I68").Value = 2160.4
$ws_LTW.Range("K68").Value = 2160.4
$ws_LTW.Range("M68").Value = -1411.4

# LTW row 71
$ws_LTW.Range("H71").Value = 2257.8572
$ws_LTW.Range("I71").Value = 2160.4
$ws_LTW.Range("K71").Value = 10802
$ws_LTW.Range("M71").Value = -7058

# LTW row 93
$ws_LTW.Range("H93").Value = 2723.2
$ws_LTW.Range("I93").Value = 2723.2
$ws_LTW.Range("K93").Value = 2723.2
$ws_LTW.Range("M93").Value = -1475.2

# LTW row 100
$ws_LTW.Range("H100").Value = 25000
$ws_LTW.Range("I100").Value = 6000
$ws_LTW.Range("J100").Value = 44000
$ws_LTW.Range("K100").Value = 6000
$ws_LTW.Range("L100").Value = 44000
$ws_LTW.Range("M100").Value = -5459
$ws_LTW.Range("N100").Value = -45082

# LTW row 132
$ws_LTW.Range("H132").Value = 5557.143
$ws_LTW.Range("I132").Value = 5566.6665
$ws_LTW.Range("K132").Value = 16699.9995
$ws_LTW.Range("M132").Value = -14169.9995

# WVR row 29
$ws_WVR.Range("H29").Value = 100
$ws_WVR.Range("J29").Value = 100
$ws_WVR.Range("L29").Value = 100
$ws_WVR.Range("N29").Value = -680

# WVR row 69
$ws_WVR.Range("H69").Value = 14249.75
$ws_WVR.Range("I69").Value = 0
$ws_WVR.Range("J69").Value = 14249.75
$ws_WVR.Range("K69").Value = 0
$ws_WVR.Range("L69").Value = 14249.75
$ws_WVR.Range("M69").ClearContents()
$ws_WVR.Range("N69").Value = -15747.75

# WVR row 72
$ws_WVR.Range("H72").Value = 14249.75
$ws_WVR.Range("I72").Value = 0
$ws_WVR.Range("J72").Value = 14249.75
$ws_WVR.Range("K72").Value = 0
$ws_WVR.Range("L72").Value = 42749.25
$ws_WVR.Range("M72").ClearContents()
$ws_WVR.Range("N72").Value = -50237.25

# WVR row 74
$ws_WVR.Range("H74").Value = 42999
$ws_WVR.Range("J74").Value = 42999
$ws_WVR.Range("L74").Value = 42999
$ws_WVR.Range("N74").Value = -44871

# WVR row 77
$ws_WVR.Range("H77").Value = 42999
$ws_WVR.Range("J77").Value = 42999
$ws_WVR.Range("L77").Value = 128997
$ws_WVR.Range("N77").Value = -138357

# WVR row 81
$ws_WVR.Range("H81").Value = 27500.5
$ws_WVR.Range("I81").Value = 32000.6
$ws_WVR.Range("J81").Value = 5000
$ws_WVR.Range("K81").Value = 64001.2
$ws_WVR.Range("L81").Value = 10000
$ws_WVR.Range("M81").Value = -62940.2
$ws_WVR.Range("N81").Value = -12122

# WVR row 84
$ws_WVR.Range("H84").Value = 27500.5
$ws_WVR.Range("I84").Value = 32000.6
$ws_WVR.Range("J84").Value = 5000
$ws_WVR.Range("K84").Value = 320006
$ws_WVR.Range("L84").Value = 50000
$ws_WVR.Range("M84").Value = -314702
$ws_WVR.Range("N84").Value = -60608

# WVR row 117
$ws_WVR.Range("H117").Value = 0
$ws_WVR.Range("J117").Value = 0
$ws_WVR.Range("L117").Value = 0
$ws_WVR.Range("N117").ClearContents()

# WVR row 126
$ws_WVR.Range("H126").Value = 3174.75
$ws_WVR.Range("I126").Value = 2233
$ws_WVR.Range("J126").Value = 6000
$ws_WVR.Range("K126").Value = 6699
$ws_WVR.Range("L126").Value = 18000
$ws_WVR.Range("M126").Value = -4229
$ws_WVR.Range("N126").Value = -22940

